# TAC-3791  Fix validation issues for import shipment excel
#
# - Bilingual (EN/AR) header labels in A1/B1
# - Drop the sample data rows (2 & 3) and clear the leftover sample row (4),
#   leaving only the formatted-but-empty B4 cell
# - Widen columns A/B so the longer bilingual headers are readable
# - Row 4 (now empty) is left selected, matching the last on-screen action

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header labels to include the Arabic translation.
$ws.Range("A1").Value = "Trip Reference*  رقم الرحلة"
$ws.Range("B1").Value = "Vas Name*   اسم الخدمة المضافة"

# Remove the sample rows entirely (rows 2 & 3), and blank out what used to
# be row 4's sample data, leaving only the formatted empty cell behind.
$ws.Range("A2:B3").ClearContents() | Out-Null
$ws.Range("A4:B4").ClearContents() | Out-Null

# Widen the columns to fit the new bilingual headers (no longer relying on
# bestFit - set explicit widths).
$ws.Columns.Item(1).ColumnWidth = 32.65
$ws.Columns.Item(2).ColumnWidth = 36.8

# Select the (now empty) row 4, mirroring the author's final selection.
$ws.Rows.Item(4).Select() | Out-Null
